$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet lists daily prices, newest first, starting at row 2.
# A new day (2025-12-10) is being prepended, so every existing data
# row (2..20) shifts down by one row (3..21); then row 2 gets the
# brand-new date with the same price values as the rest of the table.

$lastRow = 20

for ($r = $lastRow; $r -ge 2; $r--) {
    $destRow = $r + 1
    $dateVal = $ws.Cells.Item($r, 1).Value()
    $bVal = $ws.Cells.Item($r, 2).Value()
    $cVal = $ws.Cells.Item($r, 3).Value()
    $dVal = $ws.Cells.Item($r, 4).Value()

    $destDate = $ws.Cells.Item($destRow, 1)
    $destDate.NumberFormat = "@"
    $destDate.Value = $dateVal
    $destDate.ClearFormats()

    $ws.Cells.Item($destRow, 2).Value = $bVal
    $ws.Cells.Item($destRow, 3).Value = $cVal
    $ws.Cells.Item($destRow, 4).Value = $dVal
}

# New top data row: today's price snapshot.
$newDate = $ws.Cells.Item(2, 1)
$newDate.NumberFormat = "@"
$newDate.Value = "2025-12-10"
$newDate.ClearFormats()

$ws.Cells.Item(2, 2).Value = 783.5
$ws.Cells.Item(2, 3).Value = 1112
$ws.Cells.Item(2, 4).Value = 3610
